# The JSON-export search logic was factorised into a helper function; as a
# result the two admin-monitoring permission entries on the "Shared
# workspaces" sheet are now written out in swapped order (the Group entry
# that used to be on row 14 is now on row 15, and the User entry that used
# to be on row 15 is now on row 14).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Shared workspaces")

# Row 14 (was the Group row) -> becomes the User row
$ws.Range("A14").Value = "AdminInsights-6a0f4001-4816-48e9-a868-b7795aafe110"
$ws.Range("B14").Value = "AdminInsights-6a0f4001-4816-48e9-a868-b7795aafe110"
$ws.Range("C14").Value = "N/A THROUGH API"
$ws.Range("D14").Value = "User"

# Row 15 (was the User row) -> becomes the Group row
$ws.Range("A15").Value = ""
$ws.Range("B15").Value = "56c6f8c0-ded6-4444-91e6-4835e4023b53"
$ws.Range("C15").Value = "N/A THROUGH API"
$ws.Range("D15").Value = "Group"
